$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Ben Simmons -> Stephen Curry (and his position/team)
$ws.Range("A2").Value = "Stephen Curry"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Golden State Warriors"

# Row 10: Jalen Duren -> Daniel Gafford
$ws.Range("A10").Value = "Daniel Gafford"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Dallas Mavericks"

# Row 16: Stephen Curry -> Jalen Duren
$ws.Range("A16").Value = "Jalen Duren"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Detroit Pistons"

# Row 18 (last, "Daniel Gafford, PF,C, Dallas Mavericks") removed entirely.
$ws.Range("A18:C18").Delete()
